$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# "Volume 30   Number  4" -> "Volume 30   Number  5"
$ws.Range("A8").Value = "Volume 30   Number  5"
# "Report Covering the Week  1/23/2023  Through  1/29/2023" -> "...1/30/2023  Through  2/5/2023"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Crime statistics table updates (rows 15-27) ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("L15").Value = 50
$ws.Range("H16").Copy()
$ws.Range("L15").PasteSpecial(-4122)

$ws.Range("N15").Value = -25

$ws.Range("C16").Value = 7

$ws.Range("D16").Value = 8

$ws.Range("E16").Value = -12.5

$ws.Range("F16").Value = 11

$ws.Range("G16").Value = 18

$ws.Range("H16").Value = -38.888888888888

$ws.Range("I16").Value = 14

$ws.Range("J16").Value = 25

$ws.Range("K16").Value = -44

$ws.Range("L16").Value = 16.666666666666

$ws.Range("M16").Value = 27.272727272727

$ws.Range("N16").Value = -83.720930232558

$ws.Range("C17").Value = 5

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("D14").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("F17").Value = 14

$ws.Range("G17").Value = 12

$ws.Range("H17").Value = 16.666666666666

$ws.Range("I17").Value = 18

$ws.Range("K17").Value = 20

$ws.Range("L17").Value = 50

$ws.Range("M17").Value = 80

$ws.Range("N17").Value = -72.307692307692

$ws.Range("C18").Value = 3

$ws.Range("D18").Value = 2

$ws.Range("E18").Value = 50

$ws.Range("F18").Value = 20

$ws.Range("G18").Value = 20

$ws.Range("H18").Value = 0

$ws.Range("I18").Value = 31

$ws.Range("J18").Value = 33

$ws.Range("K18").Value = -6.060606060606

$ws.Range("L18").Value = 14.814814814814

$ws.Range("M18").Value = -8.823529411764

$ws.Range("N18").Value = -71.028037383177

$ws.Range("C19").Value = 15

$ws.Range("D19").Value = 16

$ws.Range("E19").Value = -6.25

$ws.Range("F19").Value = 68

$ws.Range("G19").Value = 68

$ws.Range("H19").Value = 0

$ws.Range("I19").Value = 80

$ws.Range("J19").Value = 81

$ws.Range("K19").Value = -1.234567901234

$ws.Range("L19").Value = 81.818181818181

$ws.Range("M19").Value = -1.234567901234

$ws.Range("N19").Value = -47.019867549668

$ws.Range("D20").Value = 1
$ws.Range("F16").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = -100
$ws.Range("H16").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("F20").PasteSpecial(-4122)

$ws.Range("G20").Value = 3

$ws.Range("H20").Value = -100

$ws.Range("J20").Value = 6

$ws.Range("K20").Value = -66.666666666666

$ws.Range("L20").Value = -66.666666666666

$ws.Range("N20").Value = -95.348837209302

$ws.Range("C21").Value = 30

$ws.Range("D21").Value = 27

$ws.Range("E21").Value = 11.111111111111

$ws.Range("F21").Value = 116

$ws.Range("G21").Value = 124

$ws.Range("H21").Value = -6.451612903225

$ws.Range("I21").Value = 148

$ws.Range("J21").Value = 163

$ws.Range("K21").Value = -9.202453987730

$ws.Range("L21").Value = 42.307692307692

$ws.Range("M21").Value = 5.714285714285

$ws.Range("N21").Value = -67.685589519650

$ws.Range("C22").Value = 1
$ws.Range("F16").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("F22").Value = 2

$ws.Range("H22").Value = -33.333333333333

$ws.Range("I22").Value = 2

$ws.Range("K22").Value = -33.333333333333

$ws.Range("L22").Value = 100
$ws.Range("H16").Copy()
$ws.Range("L22").PasteSpecial(-4122)

$ws.Range("M22").Value = -33.333333333333

$ws.Range("C23").Value = 2

$ws.Range("E23").Value = -33.333333333333

$ws.Range("G23").Value = 14

$ws.Range("H23").Value = -28.571428571428

$ws.Range("I23").Value = 12

$ws.Range("J23").Value = 15

$ws.Range("K23").Value = -20

$ws.Range("L23").Value = -20

$ws.Range("M23").Value = -7.692307692307

$ws.Range("C24").Value = 34

$ws.Range("D24").Value = 19

$ws.Range("E24").Value = 78.947368421052

$ws.Range("F24").Value = 109

$ws.Range("G24").Value = 115

$ws.Range("H24").Value = -5.217391304347

$ws.Range("I24").Value = 136

$ws.Range("J24").Value = 148

$ws.Range("K24").Value = -8.108108108108

$ws.Range("L24").Value = 78.947368421052

$ws.Range("M24").Value = -8.108108108108

$ws.Range("C25").Value = 4

$ws.Range("D25").Value = 8

$ws.Range("E25").Value = -50

$ws.Range("F25").Value = 29

$ws.Range("G25").Value = 36

$ws.Range("H25").Value = -19.444444444444

$ws.Range("I25").Value = 42

$ws.Range("J25").Value = 45

$ws.Range("K25").Value = -6.666666666666

$ws.Range("L25").Value = 121.052631578947

$ws.Range("M25").Value = 5

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("D14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("F26").Value = 7

$ws.Range("H26").Value = 75

$ws.Range("I26").Value = 7

$ws.Range("K26").Value = 75

$ws.Range("L26").Value = 133.333333333333

$ws.Range("C27").Value = 1
$ws.Range("F16").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = 3

$ws.Range("E27").Value = -66.666666666666

$ws.Range("G27").Value = 7

$ws.Range("H27").Value = -57.142857142857

$ws.Range("I27").Value = 4

$ws.Range("J27").Value = 9

$ws.Range("K27").Value = -55.555555555555

$ws.Range("L27").Value = 300
